$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3001
$ws.Range("I62").Value = 2501.25
$ws.Range("K62").Value = 2501.25
$ws.Range("M62").Value = -1877.25
$ws.Range("H65").Value = 3001
$ws.Range("I65").Value = 2501.25
$ws.Range("K65").Value = 12506.25
$ws.Range("M65").Value = -9386.25
$ws.Range("H76").Value = 6985.1177
$ws.Range("I76").Value = 5405.5
$ws.Range("K76").Value = 5405.5
$ws.Range("M76").Value = -5090.5
$ws.Range("H79").Value = 6985.1177
$ws.Range("I79").Value = 5405.5
$ws.Range("K79").Value = 5405.5
$ws.Range("M79").Value = -4313.5
$ws.Range("H88").Value = 4434
$ws.Range("J88").Value = 4548.4287
$ws.Range("L88").Value = 4548.4287
$ws.Range("N88").Value = -5360.4287
$ws.Range("H91").Value = 4434
$ws.Range("J91").Value = 4548.4287
$ws.Range("L91").Value = 4548.4287
$ws.Range("N91").Value = -7356.4287
$ws.Range("H96").Value = 539.9
$ws.Range("I96").Value = 716.7143
$ws.Range("J96").Value = 127.333336
$ws.Range("K96").Value = 2150.1429
$ws.Range("L96").Value = 382.000008
$ws.Range("M96").Value = -777.1428999999998
$ws.Range("N96").Value = -3128.000008
$ws.Range("H115").Value = 558.7143
$ws.Range("I115").Value = 558.7143
$ws.Range("K115").Value = 1676.1429
$ws.Range("M115").Value = -109.1428999999998
$ws.Range("H116").Value = 5795.25
$ws.Range("I116").Value = 5350.857
$ws.Range("K116").Value = 5350.857
$ws.Range("M116").Value = -1908.857
$ws.Range("H118").Value = 1486
$ws.Range("I118").Value = 2493.3333
$ws.Range("J118").Value = 881.6
$ws.Range("K118").Value = 7479.999899999999
$ws.Range("L118").Value = 2644.8
$ws.Range("M118").Value = -5822.999899999999
$ws.Range("N118").Value = -5958.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2803.5
$ws.Range("J14").Value = 5107
$ws.Range("L14").Value = 5107
$ws.Range("N14").Value = -5457
$ws.Range("H74").Value = 2251.2144
$ws.Range("I74").Value = 2077.36
$ws.Range("K74").Value = 2077.36
$ws.Range("M74").Value = -1203.36
$ws.Range("H77").Value = 2251.2144
$ws.Range("I77").Value = 2077.36
$ws.Range("K77").Value = 10386.8
$ws.Range("M77").Value = -6018.800000000001
$ws.Range("H110").Value = 5422.5713
$ws.Range("I110").Value = 4993
$ws.Range("J110").Value = 8000
$ws.Range("K110").Value = 4993
$ws.Range("L110").Value = 8000
$ws.Range("M110").Value = -2948
$ws.Range("N110").Value = -12090
$ws.Range("H132").Value = 2323.3
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1516.1852
$ws.Range("I94").Value = 1648.9412
$ws.Range("K94").Value = 1648.9412
$ws.Range("M94").Value = -1197.9412
$ws.Range("H105").Value = 2083.889
$ws.Range("I105").Value = 2083.889
$ws.Range("K105").Value = 2083.889
$ws.Range("M105").Value = -336.8890000000001
$ws.Range("H134").Value = 1298.4474
$ws.Range("I134").Value = 1293
$ws.Range("K134").Value = 3879
$ws.Range("M134").Value = -1344

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1581
$ws.Range("I31").Value = 1508.6342
$ws.Range("J31").Value = 1877.7
$ws.Range("K31").Value = 1508.6342
$ws.Range("L31").Value = 1877.7
$ws.Range("M31").Value = -1213.6342
$ws.Range("N31").Value = -2467.7
$ws.Range("H33").Value = 3397.5
$ws.Range("I33").Value = 3196.6667
$ws.Range("J33").Value = 4000
$ws.Range("K33").Value = 3196.6667
$ws.Range("L33").Value = 4000
$ws.Range("M33").Value = -2817.6667
$ws.Range("N33").Value = -4758
$ws.Range("H34").Value = 1581
$ws.Range("I34").Value = 1508.6342
$ws.Range("J34").Value = 1877.7
$ws.Range("K34").Value = 1508.6342
$ws.Range("L34").Value = 1877.7
$ws.Range("M34").Value = -1306.6342
$ws.Range("N34").Value = -2281.7
$ws.Range("H58").Value = 2885
$ws.Range("I58").Value = 2341.818
$ws.Range("J58").Value = 4080
$ws.Range("K58").Value = 2341.818
$ws.Range("L58").Value = 4080
$ws.Range("M58").Value = -2138.818
$ws.Range("N58").Value = -4486
$ws.Range("H94").Value = 11586.1
$ws.Range("I94").Value = 25558.75
$ws.Range("K94").Value = 25558.75
$ws.Range("M94").Value = -25107.75
$ws.Range("H107").Value = 12287.667
$ws.Range("I107").Value = 1194.625
$ws.Range("J107").Value = 21162.1
$ws.Range("K107").Value = 1194.625
$ws.Range("L107").Value = 21162.1
$ws.Range("M107").Value = 725.375
$ws.Range("N107").Value = -25002.1
$ws.Range("H122").Value = 396114
$ws.Range("I122").Value = 681937.75
$ws.Range("K122").Value = 2045813.25
$ws.Range("M122").Value = -2043363.25
$ws.Range("H132").Value = 1857.2778
$ws.Range("I132").Value = 1857.2778
$ws.Range("K132").Value = 5571.8334
$ws.Range("M132").Value = -3041.8334
$ws.Range("H134").Value = 2880.5312
$ws.Range("I134").Value = 1974.9667
$ws.Range("K134").Value = 5924.9001
$ws.Range("M134").Value = -3389.9001
$ws.Range("H136").Value = 2885
$ws.Range("I136").Value = 2341.818
$ws.Range("J136").Value = 4080
$ws.Range("K136").Value = 7025.454000000001
$ws.Range("L136").Value = 12240
$ws.Range("M136").Value = -4475.454000000001
$ws.Range("N136").Value = -17340

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 89.40000000000001
$ws.Range("I2").Value = 91
$ws.Range("J2").Value = 87
$ws.Range("K2").Value = 546
$ws.Range("L2").Value = 522
$ws.Range("M2").Value = -433
$ws.Range("N2").Value = -748
$ws.Range("H122").Value = 1179.6
$ws.Range("I122").Value = 1899
$ws.Range("J122").Value = 999.75
$ws.Range("K122").Value = 17091
$ws.Range("L122").Value = 8997.75
$ws.Range("M122").Value = -14641
$ws.Range("N122").Value = -13897.75
$ws.Range("H137").Value = 2569.3635
$ws.Range("I137").Value = 2286.6155
$ws.Range("J137").Value = 2977.7778
$ws.Range("K137").Value = 6859.8465
$ws.Range("L137").Value = 8933.3334
$ws.Range("M137").Value = -1759.8465
$ws.Range("N137").Value = -19133.3334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 2087
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2087
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2087
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2423
$ws.Range("H86").Value = 80000
$ws.Range("J86").Value = 80000
$ws.Range("L86").Value = 80000
$ws.Range("N86").Value = -82372
$ws.Range("H89").Value = 80000
$ws.Range("J89").Value = 80000
$ws.Range("L89").Value = 240000
$ws.Range("N89").Value = -251856
$ws.Range("H122").Value = 8418.444
$ws.Range("I122").Value = 9442
$ws.Range("K122").Value = 28326
$ws.Range("M122").Value = -25876

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1006
$ws.Range("I16").Value = 735.55554
$ws.Range("K16").Value = 735.55554
$ws.Range("M16").Value = -565.55554
$ws.Range("H22").Value = 989.0625
$ws.Range("I22").Value = 314.6
$ws.Range("K22").Value = 314.6
$ws.Range("M22").Value = -19.60000000000002
$ws.Range("H27").Value = 989.0625
$ws.Range("I27").Value = 314.6
$ws.Range("K27").Value = 314.6
$ws.Range("M27").Value = -207.6
$ws.Range("H46").Value = 999
$ws.Range("I46").Value = 999
$ws.Range("K46").Value = 999
$ws.Range("M46").Value = -811
$ws.Range("H55").Value = 1168.25
$ws.Range("I55").Value = 1249.4286
$ws.Range("K55").Value = 1249.4286
$ws.Range("M55").Value = -1076.4286
$ws.Range("H61").Value = 2344.1177
$ws.Range("I61").Value = 1521.9375
$ws.Range("K61").Value = 1521.9375
$ws.Range("M61").Value = -1319.9375
$ws.Range("H93").Value = 7483.619
$ws.Range("I93").Value = 7163.5713
$ws.Range("J93").Value = 8123.7144
$ws.Range("K93").Value = 7163.5713
$ws.Range("L93").Value = 8123.7144
$ws.Range("M93").Value = -5915.5713
$ws.Range("N93").Value = -10619.7144
$ws.Range("H113").Value = 2344.1177
$ws.Range("I113").Value = 1521.9375
$ws.Range("K113").Value = 1521.9375
$ws.Range("M113").Value = 648.0625
$ws.Range("H122").Value = 7653.846
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H130").Value = 68665
$ws.Range("J130").Value = 68665
$ws.Range("L130").Value = 68665
$ws.Range("N130").Value = -78705
$ws.Range("H132").Value = 3117.4075
$ws.Range("I132").Value = 2224.842
$ws.Range("K132").Value = 6674.526
$ws.Range("M132").Value = -4144.526

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10336
$ws.Range("I45").Value = 7968.5
$ws.Range("K45").Value = 7968.5
$ws.Range("M45").Value = -7477.5
$ws.Range("H113").Value = 1333.5
$ws.Range("I113").Value = 1173.7
$ws.Range("J113").Value = 1599.8334
$ws.Range("K113").Value = 3521.1
$ws.Range("L113").Value = 4799.5002
$ws.Range("M113").Value = -1351.1
$ws.Range("N113").Value = -9139.5002
$ws.Range("H122").Value = 2170.7693
$ws.Range("I122").Value = 2003.4445
$ws.Range("K122").Value = 6010.333500000001
$ws.Range("M122").Value = -3560.333500000001
$ws.Range("H132").Value = 1884.8182
$ws.Range("I132").Value = 1815
$ws.Range("J132").Value = 2199
$ws.Range("K132").Value = 5445
$ws.Range("L132").Value = 6597
$ws.Range("M132").Value = -2915
$ws.Range("N132").Value = -11657

